# feat: add 2022-Q4 data
#
# - Insert a new sheet "2022-Q4" between "总计" and "2021-Q2" with the
#   quarter's fund-holding detail.
# - Update the "总计" summary sheet: the former single data row (2021-Q2)
#   is pushed down to row 3, and a new row 2 is added summarising 2022-Q4.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)   # "总计" summary sheet

# Excel's default Worksheets.Add() inserts the new sheet immediately BEFORE
# the currently active sheet. The workbook opens with "2021-Q2" selected, so
# this places the new sheet right between "总计" and "2021-Q2".
$newQ = $wb.Worksheets.Add()
$newQ.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 1) "总计" sheet: add the 2022-Q4 summary row, and move the 2021-Q2 row
#    down one row (values unchanged).
# ---------------------------------------------------------------------
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)   # xlPasteFormats

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q2"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.02

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.3

# ---------------------------------------------------------------------
# 2) "2022-Q4" sheet: fund holdings detail, formatted like "总计"'s header
#    (bold / bordered / centered style).
# ---------------------------------------------------------------------
$total.Range("B1:D1").Copy()
$newQ.Range("B1:D1").PasteSpecial(-4122)
$total.Range("B1").Copy()
$newQ.Range("E1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$newQ.Range("A2:A5").PasteSpecial(-4122)

$newQ.Range("B1").Value = "基金代码"
$newQ.Range("C1").Value = "基金名称"
$newQ.Range("D1").Value = "基金规模"
$newQ.Range("E1").Value = "股票总仓位"
$newQ.Range("F1").Value = "仓位占比"
$newQ.Range("G1").Value = "持有市值(亿元)"
$newQ.Range("H1").Value = "仓位排名"

# Fund code + numeric-text columns must stay TEXT (leading zeros / trailing
# zeros need to survive, e.g. "016563", "4.00"), so force Text format before
# entering them -- otherwise Excel auto-converts numeric-looking input.
$newQ.Range("B2:B5").NumberFormat = "@"
$newQ.Range("D2:G5").NumberFormat = "@"

$newQ.Range("A2").Value = 0
$newQ.Range("B2").Value = "210009"
$newQ.Range("C2").Value = "金鹰核心资源混合"
$newQ.Range("D2").Value = "3.14"
$newQ.Range("E2").Value = "91.78"
$newQ.Range("F2").Value = "4.00"
$newQ.Range("G2").Value = "0.1256"
$newQ.Range("H2").Value = 10

$newQ.Range("A3").Value = 1
$newQ.Range("B3").Value = "162102"
$newQ.Range("C3").Value = "金鹰中小盘精选混合"
$newQ.Range("D3").Value = "3.48"
$newQ.Range("E3").Value = "76.23"
$newQ.Range("F3").Value = "3.45"
$newQ.Range("G3").Value = "0.1201"
$newQ.Range("H3").Value = 10

$newQ.Range("A4").Value = 2
$newQ.Range("B4").Value = "210002"
$newQ.Range("C4").Value = "金鹰红利价值混合A"
$newQ.Range("D4").Value = "0.96"
$newQ.Range("E4").Value = "77.22"
$newQ.Range("F4").Value = "3.81"
$newQ.Range("G4").Value = "0.0366"
$newQ.Range("H4").Value = 10

$newQ.Range("A5").Value = 3
$newQ.Range("B5").Value = "016563"
$newQ.Range("C5").Value = "金鹰红利价值混合C"
$newQ.Range("D5").Value = "0.52"
$newQ.Range("E5").Value = "77.22"
$newQ.Range("F5").Value = "3.81"
$newQ.Range("G5").Value = "0.0198"
$newQ.Range("H5").Value = 10

# Restore the original tab selection (the "2021-Q2" sheet was the active /
# selected tab before this edit). Look it up fresh by name -- the sheet
# collection shifted position when the new sheet was inserted, so a
# reference captured before the Add() would now point at the wrong tab.
$oldQ = $wb.Worksheets.Item("2021-Q2")
$oldQ.Activate()

Write-Output "2022-Q4 sheet added"
